$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.750601
$ws.Range("H2").Value = 65.251803
$ws.Range("I2").Value = 0.3612485837751334
$ws.Range("J2").Value = 0.3773020033645041
$ws.Range("M2").Value = 0.5351276666666666
$ws.Range("N2").Value = 1.605383
$ws.Range("O2").Value = 0.04512696048487568
$ws.Range("P2").Value = 0.0511377912725271
$ws.Range("Q2").Value = 11.63934836172766
$ws.Range("R2").Value = 104.754135255549
$ws.Range("S2").Value = 0.01630205056523774
$ws.Range("T2").Value = 0.01929439109476033
$ws.Range("G3").Value = 21.750601
$ws.Range("H3").Value = 65.251803
$ws.Range("I3").Value = 0.3612485837751334
$ws.Range("J3").Value = 0.3773020033645041
$ws.Range("O3").Value = 0.1574346892487883
$ws.Range("P3").Value = 0.178404709542934
$ws.Range("Q3").Value = 40.606262258703
$ws.Range("R3").Value = 365.456360328327
$ws.Range("S3").Value = 0.05687305852820301
$ws.Range("T3").Value = 0.06731245432021146
$ws.Range("G4").Value = 21.750601
$ws.Range("H4").Value = 65.251803
$ws.Range("I4").Value = 0.3612485837751334
$ws.Range("J4").Value = 0.3773020033645041
$ws.Range("M4").Value = 2.563637666666667
$ws.Range("N4").Value = 7.690913
$ws.Range("O4").Value = 0.2161898606398702
$ws.Range("P4").Value = 0.2449859651492294
$ws.Range("Q4").Value = 55.76065999623766
$ws.Range("R4").Value = 501.845939966139
$ws.Range("S4").Value = 0.07809828098269656
$ws.Range("T4").Value = 0.09243369544699082
$ws.Range("G5").Value = 21.750601
$ws.Range("H5").Value = 65.251803
$ws.Range("I5").Value = 0.3612485837751334
$ws.Range("J5").Value = 0.3773020033645041
$ws.Range("M5").Value = 4.181529
$ws.Range("N5").Value = 8.363058000000001
$ws.Range("O5").Value = 0.3526255615314758
$ws.Range("P5").Value = 0.2663964389831199
$ws.Range("Q5").Value = 90.950768848929
$ws.Range("R5").Value = 545.704613093574
$ws.Range("S5").Value = 0.1273854847061568
$ws.Range("T5").Value = 0.100511910117501
$ws.Range("G6").Value = 21.750601
$ws.Range("H6").Value = 65.251803
$ws.Range("I6").Value = 0.3612485837751334
$ws.Range("J6").Value = 0.3773020033645041
$ws.Range("M6").Value = 2.711072333333334
$ws.Range("N6").Value = 8.133217
$ws.Range("O6").Value = 0.2286229280949899
$ws.Range("P6").Value = 0.2590750950521895
$ws.Range("Q6").Value = 58.96745260447234
$ws.Range("R6").Value = 530.707073440251
$ws.Range("S6").Value = 0.08258970899283927
$ws.Range("T6").Value = 0.09774955238504042
$ws.Range("I7").Value = 0.2797740820980411
$ws.Range("J7").Value = 0.2922068802649305
$ws.Range("M7").Value = 0.5351276666666666
$ws.Range("N7").Value = 1.605383
$ws.Range("O7").Value = 0.04512696048487568
$ws.Range("P7").Value = 0.0511377912725271
$ws.Range("Q7").Value = 9.014258187787666
$ws.Range("R7").Value = 81.12832369008899
$ws.Range("S7").Value = 0.01262535394753066
$ws.Range("T7").Value = 0.01494281445138433
$ws.Range("I8").Value = 0.2797740820980411
$ws.Range("J8").Value = 0.2922068802649305
$ws.Range("O8").Value = 0.1574346892487883
$ws.Range("P8").Value = 0.178404709542934
$ws.Range("S8").Value = 0.04404614567497009
$ws.Range("T8").Value = 0.05213108360011182
$ws.Range("I9").Value = 0.2797740820980411
$ws.Range("J9").Value = 0.2922068802649305
$ws.Range("M9").Value = 2.563637666666667
$ws.Range("N9").Value = 7.690913
$ws.Range("O9").Value = 0.2161898606398702
$ws.Range("P9").Value = 0.2449859651492294
$ws.Range("Q9").Value = 43.18463287689767
$ws.Range("R9").Value = 388.661695892079
$ws.Range("S9").Value = 0.06048431981942309
$ws.Range("T9").Value = 0.07158658458494929
$ws.Range("I10").Value = 0.2797740820980411
$ws.Range("J10").Value = 0.2922068802649305
$ws.Range("M10").Value = 4.181529
$ws.Range("N10").Value = 8.363058000000001
$ws.Range("O10").Value = 0.3526255615314758
$ws.Range("P10").Value = 0.2663964389831199
$ws.Range("Q10").Value = 70.438111078269
$ws.Range("R10").Value = 422.6286664696141
$ws.Range("S10").Value = 0.09865549280177494
$ws.Range("T10").Value = 0.07784287234894438
$ws.Range("I11").Value = 0.2797740820980411
$ws.Range("J11").Value = 0.2922068802649305
$ws.Range("M11").Value = 2.711072333333334
$ws.Range("N11").Value = 8.133217
$ws.Range("O11").Value = 0.2286229280949899
$ws.Range("P11").Value = 0.2590750950521895
$ws.Range("Q11").Value = 45.66817883041234
$ws.Range("R11").Value = 411.013609473711
$ws.Range("S11").Value = 0.06396276985434224
$ws.Range("T11").Value = 0.07570352527954062
$ws.Range("G12").Value = 8.938416999999999
$ws.Range("H12").Value = 26.815251
$ws.Range("I12").Value = 0.1484552303838214
$ws.Range("J12").Value = 0.1550523887136425
$ws.Range("M12").Value = 0.5351276666666666
$ws.Range("N12").Value = 1.605383
$ws.Range("O12").Value = 0.04512696048487568
$ws.Range("P12").Value = 0.0511377912725271
$ws.Range("Q12").Value = 4.783194232903665
$ws.Range("R12").Value = 43.048748096133
$ws.Range("S12").Value = 0.006699333315303824
$ws.Range("T12").Value = 0.007929036690344985
$ws.Range("G13").Value = 8.938416999999999
$ws.Range("H13").Value = 26.815251
$ws.Range("I13").Value = 0.1484552303838214
$ws.Range("J13").Value = 0.1550523887136425
$ws.Range("O13").Value = 0.1574346892487883
$ws.Range("P13").Value = 0.178404709542934
$ws.Range("Q13").Value = 16.687157512551
$ws.Range("R13").Value = 150.184417612959
$ws.Range("S13").Value = 0.02337200306283421
$ws.Range("T13").Value = 0.02766207637239548
$ws.Range("G14").Value = 8.938416999999999
$ws.Range("H14").Value = 26.815251
$ws.Range("I14").Value = 0.1484552303838214
$ws.Range("J14").Value = 0.1550523887136425
$ws.Range("M14").Value = 2.563637666666667
$ws.Range("N14").Value = 7.690913
$ws.Range("O14").Value = 0.2161898606398702
$ws.Range("P14").Value = 0.2449859651492294
$ws.Range("Q14").Value = 22.91486250157367
$ws.Range("R14").Value = 206.233762514163
$ws.Range("S14").Value = 0.03209451556793817
$ws.Range("T14").Value = 0.03798565909770518
$ws.Range("G15").Value = 8.938416999999999
$ws.Range("H15").Value = 26.815251
$ws.Range("I15").Value = 0.1484552303838214
$ws.Range("J15").Value = 0.1550523887136425
$ws.Range("M15").Value = 4.181529
$ws.Range("N15").Value = 8.363058000000001
$ws.Range("O15").Value = 0.3526255615314758
$ws.Range("P15").Value = 0.2663964389831199
$ws.Range("Q15").Value = 37.376249899593
$ws.Range("R15").Value = 224.257499397558
$ws.Range("S15").Value = 0.05234910897637963
$ws.Range("T15").Value = 0.04130540420914085
$ws.Range("G16").Value = 8.938416999999999
$ws.Range("H16").Value = 26.815251
$ws.Range("I16").Value = 0.1484552303838214
$ws.Range("J16").Value = 0.1550523887136425
$ws.Range("M16").Value = 2.711072333333334
$ws.Range("N16").Value = 8.133217
$ws.Range("O16").Value = 0.2286229280949899
$ws.Range("P16").Value = 0.2590750950521895
$ws.Range("Q16").Value = 24.23269503249633
$ws.Range("R16").Value = 218.094255292467
$ws.Range("S16").Value = 0.03394026946136557
$ws.Range("T16").Value = 0.04017021234405595
$ws.Range("G17").Value = 7.6853705
$ws.Range("H17").Value = 15.370741
$ws.Range("I17").Value = 0.1276437928732263
$ws.Range("J17").Value = 0.08887741190073968
$ws.Range("M17").Value = 0.5351276666666666
$ws.Range("N17").Value = 1.605383
$ws.Range("O17").Value = 0.04512696048487568
$ws.Range("P17").Value = 0.0511377912725271
$ws.Range("Q17").Value = 4.112654383133832
$ws.Range("R17").Value = 24.675926298803
$ws.Range("S17").Value = 0.00576017639712974
$ws.Range("T17").Value = 0.004544994538622442
$ws.Range("G18").Value = 7.6853705
$ws.Range("H18").Value = 15.370741
$ws.Range("I18").Value = 0.1276437928732263
$ws.Range("J18").Value = 0.08887741190073968
$ws.Range("O18").Value = 0.1574346892487883
$ws.Range("P18").Value = 0.178404709542934
$ws.Range("Q18").Value = 14.3478412425615
$ws.Range("R18").Value = 86.087047455369
$ws.Range("S18").Value = 0.02009556086553309
$ws.Range("T18").Value = 0.01585614885507917
$ws.Range("G19").Value = 7.6853705
$ws.Range("H19").Value = 15.370741
$ws.Range("I19").Value = 0.1276437928732263
$ws.Range("J19").Value = 0.08887741190073968
$ws.Range("M19").Value = 2.563637666666667
$ws.Range("N19").Value = 7.690913
$ws.Range("O19").Value = 0.2161898606398702
$ws.Range("P19").Value = 0.2449859651492294
$ws.Range("Q19").Value = 19.70250529608883
$ws.Range("R19").Value = 118.215031776533
$ws.Range("S19").Value = 0.02759529379280725
$ws.Range("T19").Value = 0.02177371853446831
$ws.Range("G20").Value = 7.6853705
$ws.Range("H20").Value = 15.370741
$ws.Range("I20").Value = 0.1276437928732263
$ws.Range("J20").Value = 0.08887741190073968
$ws.Range("M20").Value = 4.181529
$ws.Range("N20").Value = 8.363058000000001
$ws.Range("O20").Value = 0.3526255615314758
$ws.Range("P20").Value = 0.2663964389831199
$ws.Range("Q20").Value = 32.1365996214945
$ws.Range("R20").Value = 128.546398485978
$ws.Range("S20").Value = 0.04501046413792882
$ws.Range("T20").Value = 0.02367662603639302
$ws.Range("G21").Value = 7.6853705
$ws.Range("H21").Value = 15.370741
$ws.Range("I21").Value = 0.1276437928732263
$ws.Range("J21").Value = 0.08887741190073968
$ws.Range("M21").Value = 2.711072333333334
$ws.Range("N21").Value = 8.133217
$ws.Range("O21").Value = 0.2286229280949899
$ws.Range("P21").Value = 0.2590750950521895
$ws.Range("Q21").Value = 20.83559533396617
$ws.Range("R21").Value = 125.013572003797
$ws.Range("S21").Value = 0.02918229767982741
$ws.Range("T21").Value = 0.02302592393617673
$ws.Range("G22").Value = 4.990062666666667
$ws.Range("H22").Value = 14.970188
$ws.Range("I22").Value = 0.08287831086977776
$ws.Range("J22").Value = 0.08656131575618316
$ws.Range("M22").Value = 0.5351276666666666
$ws.Range("N22").Value = 1.605383
$ws.Range("O22").Value = 0.04512696048487568
$ws.Range("P22").Value = 0.0511377912725271
$ws.Range("Q22").Value = 2.670320591333777
$ws.Range("R22").Value = 24.032885322004
$ws.Range("S22").Value = 0.003740046259673703
$ws.Range("T22").Value = 0.004426554497415006
$ws.Range("G23").Value = 4.990062666666667
$ws.Range("H23").Value = 14.970188
$ws.Range("I23").Value = 0.08287831086977776
$ws.Range("J23").Value = 0.08656131575618316
$ws.Range("O23").Value = 0.1574346892487883
$ws.Range("P23").Value = 0.178404709542934
$ws.Range("Q23").Value = 9.315962962587999
$ws.Range("R23").Value = 83.843666663292
$ws.Range("S23").Value = 0.01304792111724794
$ws.Range("T23").Value = 0.01544294639513605
$ws.Range("G24").Value = 4.990062666666667
$ws.Range("H24").Value = 14.970188
$ws.Range("I24").Value = 0.08287831086977776
$ws.Range("J24").Value = 0.08656131575618316
$ws.Range("M24").Value = 2.563637666666667
$ws.Range("N24").Value = 7.690913
$ws.Range("O24").Value = 0.2161898606398702
$ws.Range("P24").Value = 0.2449859651492294
$ws.Range("Q24").Value = 12.79271261129378
$ws.Range("R24").Value = 115.134413501644
$ws.Range("S24").Value = 0.01791745047700509
$ws.Range("T24").Value = 0.02120630748511573
$ws.Range("G25").Value = 4.990062666666667
$ws.Range("H25").Value = 14.970188
$ws.Range("I25").Value = 0.08287831086977776
$ws.Range("J25").Value = 0.08656131575618316
$ws.Range("M25").Value = 4.181529
$ws.Range("N25").Value = 8.363058000000001
$ws.Range("O25").Value = 0.3526255615314758
$ws.Range("P25").Value = 0.2663964389831199
$ws.Range("Q25").Value = 20.866091752484
$ws.Range("R25").Value = 125.196550514904
$ws.Range("S25").Value = 0.0292250109092356
$ws.Range("T25").Value = 0.02305962627114062
$ws.Range("G26").Value = 4.990062666666667
$ws.Range("H26").Value = 14.970188
$ws.Range("I26").Value = 0.08287831086977776
$ws.Range("J26").Value = 0.08656131575618316
$ws.Range("M26").Value = 2.711072333333334
$ws.Range("N26").Value = 8.133217
$ws.Range("O26").Value = 0.2286229280949899
$ws.Range("P26").Value = 0.2590750950521895
$ws.Range("Q26").Value = 13.52842083719956
$ws.Range("R26").Value = 121.755787534796
$ws.Range("S26").Value = 0.01894788210661542
$ws.Range("T26").Value = 0.02242588110737574
